$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.317.66'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.551.22'
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '209.71'
$ws.Range('E5').Value = '  -1.54%  '
$ws.Range('E6').Value = '  -1.85%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.68'
$ws.Range('E8').Value = '  -2.10%  '
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('E10').Value = '  -1.51%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0889'
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.773.17'
$ws.Range('E12').Value = '  -1.36%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.555.05'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.302.45'
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('E16').Value = '  -2.47%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '60.55'
$ws.Range('E17').Value = '  -3.00%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '227.96'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.33'
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.92'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.88'
$ws.Range('E23').Value = '  -3.15%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.03'
$ws.Range('E24').Value = '  -4.90%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '150.60'
$ws.Range('E25').Value = '  -0.82%  '
$ws.Range('E26').Value = '  -1.97%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.25'
$ws.Range('E29').Value = '  -3.15%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0467'
$ws.Range('E30').Value = '  -4.05%  '
$ws.Range('E31').Value = '  -4.59%  '
$ws.Range('E32').Value = '  -1.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.03'
$ws.Range('E33').Value = '  -2.66%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.385.67'
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('E36').Value = '  -3.27%  '
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('E39').Value = '  -3.11%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.93'
$ws.Range('E40').Value = '  +1.03%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.510'
$ws.Range('E41').Value = '  -3.07%  '
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.777'
$ws.Range('E43').Value = '  -2.14%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0466'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.40'
$ws.Range('E45').Value = '  -2.28%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '61.91'
$ws.Range('E46').Value = '  -2.12%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.686.85'
$ws.Range('E47').Value = '  -1.46%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.906'
$ws.Range('E48').Value = '  -6.22%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '85.69'
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₆0103'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '41.89'
$ws.Range('E51').Value = '  +4.90%  '
